# Atualização de bases das ligas, do dia: 11-04-2024 às 00:31
#
# Three pairs of rows had their match records swapped (everything except
# the running index in column A), and a handful of isolated odds cells
# were corrected on six other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2) {
    $r1 = $ws.Range("B$row1`:AC$row1")
    $r2 = $ws.Range("B$row2`:AC$row2")
    $v1 = $r1.Value2
    $v2 = $r2.Value2
    $r1.Value = $v2
    $r2.Value = $v1
}

# Row 159 <-> Row 160 (match ids 6810145 / 6810142)
Swap-Rows 159 160

# Row 181 <-> Row 182 (match ids 6810167 / 6810169)
Swap-Rows 181 182

# Row 183 <-> Row 184 (match ids 6810165 / 6810168)
Swap-Rows 183 184

# Isolated odds corrections on rows 259, 260, 261, 264, 265, 266
$ws.Range("R259").Value = 1.8
$ws.Range("S259").Value = 2.05

$ws.Range("U260").Value = 1.95
$ws.Range("V260").Value = 1.9

$ws.Range("N261").Value = 1.727
$ws.Range("P261").Value = 4.75
$ws.Range("R261").Value = 1.975
$ws.Range("S261").Value = 1.875
$ws.Range("U261").Value = 1.975
$ws.Range("V261").Value = 1.875

$ws.Range("N264").Value = 3
$ws.Range("O264").Value = 3.2
$ws.Range("U264").Value = 2.05
$ws.Range("V264").Value = 1.8

$ws.Range("U265").Value = 1.875
$ws.Range("V265").Value = 1.975

$ws.Range("R266").Value = 1.975
$ws.Range("S266").Value = 1.875
